# Insert a new record (most recent date) at the top of the Achicoria -
# Vega Modelo de Temuco price history block. All the existing rows from
# 87 downward shift down by one row (handled automatically by Insert()).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 87, pushing rows 87:146 down to 88:147
$ws.Rows.Item(87).Insert()

# Populate the newly inserted row 87 with the new record's data
$ws.Cells.Item(87, 1).Value = 10
$ws.Cells.Item(87, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(87, 3).Value = "La Araucanía"
$ws.Cells.Item(87, 4).Value = 45176
$ws.Cells.Item(87, 5).Value = 9
$ws.Cells.Item(87, 6).Value = 100112010
$ws.Cells.Item(87, 7).Value = "Achicoria"
$ws.Cells.Item(87, 8).Value = "Sin especificar"
$ws.Cells.Item(87, 9).Value = "Primera"
$ws.Cells.Item(87, 10).Value = 300
$ws.Cells.Item(87, 11).Value = 10000
$ws.Cells.Item(87, 12).Value = 10000
$ws.Cells.Item(87, 13).Value = 10000
$ws.Cells.Item(87, 14).Value = "`$/caja 18 unidades"
$ws.Cells.Item(87, 15).Value = "Región Metropolitana"
$ws.Cells.Item(87, 16).Value = 556
$ws.Cells.Item(87, 17).Value = 18
$ws.Cells.Item(87, 18).Value = "Hortaliza"
